# Manual-extraction update: append newly-extracted lobby rows to the
# per-session sheets and to the aggregated "Lobbies_Trustworthy_Lobbies"
# sheet, matching additional rows found while re-checking the SRT files.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 2022-02-12_S1 (sheet11) -- add rows 8 and 9
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-02-12_S1")
$ws.Range("H8").Value = 15
$ws.Range("I8:J9").NumberFormat = "h:mm:ss"
$ws.Range("I8").Value = 0.1006712962962963
$ws.Range("J8").Value = 0.10765046296296295
$ws.Range("H9").Value = 20
$ws.Range("I9").Value = 0.13614583333333333
$ws.Range("J9").Value = 0.14758101851851851

# ---------------------------------------------------------------------
# 2022-02-19_S1 (sheet15) -- add row 9
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-02-19_S1")
$ws.Range("H9").Value = 7
$ws.Range("I9:J9").NumberFormat = "h:mm:ss"
$ws.Range("I9").Value = 0.04608796296296296
$ws.Range("J9").Value = 0.05145833333333333

# ---------------------------------------------------------------------
# 2022-02-21_S1 (sheet16) -- add row 13, clear/empty C3 style, fix J6
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-02-21_S1")
$ws.Range("C3").NumberFormat = "h:mm"
$ws.Range("J6").Value = 0.0026041666666666665
$ws.Range("H13").Value = 13
$ws.Range("I13:J13").NumberFormat = "h:mm:ss"
$ws.Range("I13").Value = 0.05693287037037037
$ws.Range("J13").Value = 0.0628125

# ---------------------------------------------------------------------
# 2022-02-23_S1 (sheet18) -- add rows 6 and 7
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-02-23_S1")
$ws.Range("H6").Value = 11
$ws.Range("I6:J7").NumberFormat = "h:mm:ss"
$ws.Range("I6").Value = 0.06386574074074074
$ws.Range("J6").Value = 0.06899305555555556
$ws.Range("H7").Value = 15
$ws.Range("I7").Value = 0.09153935185185186
$ws.Range("J7").Value = 0.09533564814814816

# ---------------------------------------------------------------------
# 2022-03.01_S1 (sheet27) -- add row 10
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-03.01_S1")
$ws.Range("H10").Value = 17
$ws.Range("I10:J10").NumberFormat = "h:mm:ss"
$ws.Range("I10").Value = 0.10886574074074074
$ws.Range("J10").Value = 0.11337962962962962

# ---------------------------------------------------------------------
# 2022-03-02_S1 (sheet28) -- add row 7
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-03-02_S1")
$ws.Range("H7").Value = 12
$ws.Range("I7:J7").NumberFormat = "h:mm:ss"
$ws.Range("I7").Value = 0.09158564814814814
$ws.Range("J7").Value = 0.09711805555555557

# ---------------------------------------------------------------------
# 2022-05-19_S1 (sheet32) -- add row 12
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-05-19_S1")
$ws.Range("H12").Value = 18
$ws.Range("I12:J12").NumberFormat = "h:mm:ss"
$ws.Range("I12").Value = 0.09965277777777777
$ws.Range("J12").Value = 0.10324074074074074

# ---------------------------------------------------------------------
# Lobbies_Trustworthy_Lobbies (sheet26) -- append the same newly found
# rows (rows 171-179), one block per source session above.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Lobbies_Trustworthy_Lobbies")
$ws.Range("C171:D179").NumberFormat = "h:mm:ss"

$ws.Range("A171").Value = "2022-02-23_S1_vikramafc_1307860812"
$ws.Range("B171").Value = 11
$ws.Range("C171").Value = 0.06386574074074074
$ws.Range("D171").Value = 0.06899305555555556

$ws.Range("A172").Value = "2022-02-23_S1_vikramafc_1307860812"
$ws.Range("B172").Value = 15
$ws.Range("C172").Value = 0.09153935185185186
$ws.Range("D172").Value = 0.09533564814814816

$ws.Range("A173").Value = "2022-02-21_S1_aribunnie_1305616488"
$ws.Range("B173").Value = 13
$ws.Range("C173").Value = 0.05693287037037037
$ws.Range("D173").Value = 0.0628125

$ws.Range("A174").Value = "2022-03-01_S1_aribunnie_1412358486"
$ws.Range("B174").Value = 17
$ws.Range("C174").Value = 0.10886574074074074
$ws.Range("D174").Value = 0.11337962962962962

$ws.Range("A175").Value = "2022-02-12_S1_ressnie_1295839771"
$ws.Range("B175").Value = 15
$ws.Range("C175").Value = 0.1006712962962963
$ws.Range("D175").Value = 0.10765046296296295

$ws.Range("A176").Value = "2022-02-12_S1_ressnie_1295839771"
$ws.Range("B176").Value = 20
$ws.Range("C176").Value = 0.13614583333333333
$ws.Range("D176").Value = 0.14758101851851851

$ws.Range("A177").Value = "2022-05-19_S1_br00d_1488152742"
$ws.Range("B177").Value = 18
$ws.Range("C177").Value = 0.09965277777777777
$ws.Range("D177").Value = 0.10324074074074074

$ws.Range("A178").Value = "2022-02-19_S1_willyutv_1303482578"
$ws.Range("B178").Value = 7
$ws.Range("C178").Value = 0.04608796296296296
$ws.Range("D178").Value = 0.05145833333333333

$ws.Range("A179").Value = "2022-03-02_S1_irepptar_1413360539"
$ws.Range("B179").Value = 12
$ws.Range("C179").Value = 0.09158564814814814
$ws.Range("D179").Value = 0.09711805555555557

# ---------------------------------------------------------------------
# View / selection bookkeeping to mirror the author's final state:
# each touched session sheet ends with its header cell selected, the
# summary sheet lands on the newly-appended row, and the workbook's
# final active sheet is 2022-02-21_S1 (was Main_Streamer_Lobbies).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("2022-02-12_S1")
$ws.Range("G2").Select()

$ws = $wb.Worksheets.Item("2022-02-19_S1")
$ws.Range("G2").Select()

$ws = $wb.Worksheets.Item("2022-02-23_S1")
$ws.Range("G2").Select()

$ws = $wb.Worksheets.Item("2022-03.01_S1")
$ws.Range("B2").Select()

$ws = $wb.Worksheets.Item("2022-03-02_S1")
$ws.Range("G2").Select()

$ws = $wb.Worksheets.Item("2022-05-19_S1")
$ws.Range("G2").Select()

$ws = $wb.Worksheets.Item("Lobbies_Trustworthy_Lobbies")
$ws.Range("A179").Select()

$ws = $wb.Worksheets.Item("2022-02-21_S1")
$ws.Activate()
$ws.Range("I6:J6").Select()
